$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Save" column (H) to the s_vals sheet, mirroring the header
# formatting used by the other header cells (e.g. G1's bold/bordered style).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Fill in the Save values for each data row.
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
